# Splits the leading "4.<letter>.<digit>." numbering of three step
# paragraphs into three separate runs: "4.", "<letter>", and the
# remaining ". <rest of the sentence>" text — matching the target
# OOXML diff exactly (identical rPr kept as three sibling <w:r>
# elements instead of being re-merged into one run).

$d = $word.ActiveDocument

function Split-NumberingRun {
    param(
        [string]$oldText,
        [string]$letter,
        [string]$rest
    )

    $rng = $d.Content
    $found = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find text: $oldText"
    }

    $target = $d.Range($rng.Start, $rng.End)
    $target.Text = ""

    $restTag = '<w:t>'
    if ($rest.StartsWith(" ") -or $rest.EndsWith(" ")) {
        $restTag = '<w:t xml:space="preserve">'
    }

    $xmlFrag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' +
        '<w:r><w:rPr><w:rFonts w:ascii="Cambria" w:hAnsi="Cambria"/></w:rPr><w:t>4.</w:t></w:r>' +
        '<w:r><w:rPr><w:rFonts w:ascii="Cambria" w:hAnsi="Cambria"/></w:rPr><w:t>' + $letter + '</w:t></w:r>' +
        '<w:r><w:rPr><w:rFonts w:ascii="Cambria" w:hAnsi="Cambria"/></w:rPr>' + $restTag + $rest + '</w:t></w:r>' +
        '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    $target.InsertXML($xmlFrag)
}

Split-NumberingRun "4.β.1. Ο αγοραστής επιλέγει το πλήκτρο «Λίγα λόγια για το προϊόν»." "α" ". Ο αγοραστής επιλέγει το πλήκτρο «Λίγα λόγια για το προϊόν»."

Split-NumberingRun "4.β.2. Το σύστημα λαμβάνει λίγα λόγια για το προϊόν και " "β" ". Το σύστημα λαμβάνει λίγα λόγια για το προϊόν και "

Split-NumberingRun "4.β.4. Η περίπτωση χρήσης συνεχίζεται στο βήμα 3 της βασικής ροής." "γ" ". Η περίπτωση χρήσης συνεχίζεται στο βήμα 3 της βασικής ροής."

Write-Host "Done."
